# Update tab names in all BOMs, fix bi-color LED naming.
# (This workbook's slice of that commit: rename the sheet tab and
#  refresh the saved view/selection state.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "KickAid" to "BOM"
$ws.Name = "BOM"

# Make sure it's the active/selected sheet
$ws.Activate()

# Move the view so row 39 scrolls to the top of the window, then land
# the active selection on D67 (mirrors the saved sheetView/selection
# state: topLeftCell="A39", activeCell="D67").
$excel.Goto($ws.Range("D67"), $true)

$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1

$ws.Range("D67").Select()
